# Applies the "Sincronização de dados e código" update:
#  - quotations!3 -> quotation 1774f003-3175-4be8-801c-c7e958af19a8 gets approved
#  - items sheet rows 3-9 get re-synced (new item ids + a couple of
#    productService references rotated between rows, with their
#    quantity/value/type following the reference)
#
# Note: the id-like columns (item_id / productService_id on the items
# sheet, and the various *.id columns elsewhere) are stored as
# base64-encoded "<uuid>:<accountId>" strings, so the literals below are
# base64 already, matching the source data's own encoding.

$wb = $excel.ActiveWorkbook

$wsQuotations = $wb.Worksheets.Item("quotations")
$wsItems      = $wb.Worksheets.Item("items")

# ---------------------------------------------------------------------
# 1. quotations sheet, row 3 (quotation id 1774f003-3175-4be8-801c-c7e958af19a8)
# ---------------------------------------------------------------------
$wsQuotations.Range("I3").Value = "Aprovada"
$wsQuotations.Range("L3").Value = "Edina K O Mishima"
$wsQuotations.Range("N3").Value = "accounts/57016/quotations/1774f003-3175-4be8-801c-c7e958af19a8/signatures/650d7c22-e91c-4757-9ac9-9ad5beab0def.png"
$wsQuotations.Range("U3").Value = "approved"

# ---------------------------------------------------------------------
# 2. items sheet, rows 3-9 re-sync
# ---------------------------------------------------------------------

# Row 3: only the item_id changes (fresh id from the source system)
$wsItems.Range("A3").Value = "MDU3Zjg5ODMtM2E3Mi00NDM4LWIyODctOGNkNDhiZTkyZjU5OjU3MDE2"

# Row 4: new item_id; productService reference rotates in from row 6
$wsItems.Range("A4").Value = "MDgwZWUyMTktODkzYi00ZTlmLTkwYmEtMjkzOTdiZDA0YWM1OjU3MDE2"
$wsItems.Range("B4").Value = 1
$wsItems.Range("C4").Value = 9
$wsItems.Range("G4").Value = "YjEzMmNlNWQtOGU2Ny00NmIwLWJiYzMtNDk0Zjg1YzMyNGIyOjU3MDE2"
$wsItems.Range("H4").Value = 9

# Row 5: new item_id; productService reference rotates in from row 4
$wsItems.Range("A5").Value = "MGE0MDEzNzktZGQyZi00NTllLWIwMmMtOWVmZDEwODUyZmU3OjU3MDE2"
$wsItems.Range("B5").Value = 2
$wsItems.Range("C5").Value = 30
$wsItems.Range("G5").Value = "YmVjYmU1OWUtZDQ0Zi00ZWNlLTgxZWMtNDA5MDk3NjNlMGM0OjU3MDE2"
$wsItems.Range("H5").Value = 15

# Row 6: new item_id; productService reference rotates in from row 5
$wsItems.Range("A6").Value = "NmI4ZjZmNGQtODZhYi00YmRkLTk5NWEtNjk3YjdkOTA4YmRkOjU3MDE2"
$wsItems.Range("C6").Value = 25200000000000000
$wsItems.Range("G6").Value = "NDRhZmYxYjItOTQxZi00NTQ1LTk2M2MtMzEzYzY2ODdjODgyOjU3MDE2"
$wsItems.Range("H6").Value = 25200000000000000

# Row 7: new item_id; productService reference swaps in from row 8
$wsItems.Range("A7").Value = "OWVmOTBjMjAtMDY3Zi00YWEyLWFjMjEtNGE3YmRkMjY3ZWQ5OjU3MDE2"
$wsItems.Range("C7").Value = 204
$wsItems.Range("G7").Value = "ZTg2MzY3YjMtNGQ4Ni00MjViLWJhODQtNzk5OWM5NDgwY2Q3OjU3MDE2"
$wsItems.Range("H7").Value = 204
$wsItems.Range("I7").Value = "product"

# Row 8: new item_id; productService reference swaps in from row 7
$wsItems.Range("A8").Value = "ZDU2NGE2M2EtYmE1OS00ZjRmLThhNDgtMjIxMjYzNzliNDZhOjU3MDE2"
$wsItems.Range("C8").Value = 350
$wsItems.Range("G8").Value = "YmJmNjhlOWMtMWYyMS00MTVlLWI5MzctN2NmMDNhNDY5OWFkOjU3MDE2"
$wsItems.Range("H8").Value = 350
$wsItems.Range("I8").Value = "service"

# Row 9: only the item_id changes (fresh id from the source system)
$wsItems.Range("A9").Value = "ZTZlYTVkODUtZTkwMy00MWE0LWJiNmYtYTM0NGM4NGM2YTkyOjU3MDE2"
